# Update "database and change read_price algorithm"
#
# The yearly table reports 5 trailing twelve-month (TTM) periods per line,
# in columns E:I (oldest -> newest). The update rolls the window forward by
# one fiscal year: the oldest period (1396/12) is dropped, the remaining
# four periods shift one column to the left (F->E, G->F, H->G, I->H), and a
# brand-new trailing period (1401/12) is appended in column I with freshly
# read values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rows: shift the "twelve months ended .../12" period labels ---
$headerRows = @(8, 24)
foreach ($r in $headerRows) {
    $ws.Cells.Item($r, 5).Value = "دوازده ماهه منتهی به 1397/12"   # E
    $ws.Cells.Item($r, 6).Value = "دوازده ماهه منتهی به 1398/12"   # F
    $ws.Cells.Item($r, 7).Value = "دوازده ماهه منتهی به 1399/12"   # G
    $ws.Cells.Item($r, 8).Value = "دوازده ماهه منتهی به 1400/12"   # H
    $ws.Cells.Item($r, 9).Value = "دوازده ماهه منتهی به 1401/12"   # I
}

# --- Data rows: new values for columns E:I (oldest dropped, newest read) ---
$newRowValues = @{
    10 = @(0,      0,      26317,  48508,   87760)
    11 = @(0,      0,      0,      0,       0)
    12 = @(0,      0,      0,      0,       0)
    13 = @(0,      0,      0,      0,       0)
    14 = @(0,      0,      61764,  139139,  329476)
    15 = @(0,      0,      0,      0,       0)
    16 = @(111157, 115440, 61145,  65590,   89704)
    17 = @(125875, 179776, 362831, 718123,  1104707)
    18 = @(0,      0,      0,      0,       0)
    19 = @(390266, 463864, 497441, 1221585, 1624139)
    20 = @(627298, 759080, 1009498,2192945, 3235786)
    26 = @(735,    768,    722,    852,     644)
    27 = @(395,    416,    416,    520,     788)
}

foreach ($r in $newRowValues.Keys) {
    $vals = $newRowValues[$r]
    $ws.Cells.Item($r, 5).Value = $vals[0]   # E
    $ws.Cells.Item($r, 6).Value = $vals[1]   # F
    $ws.Cells.Item($r, 7).Value = $vals[2]   # G
    $ws.Cells.Item($r, 8).Value = $vals[3]   # H
    $ws.Cells.Item($r, 9).Value = $vals[4]   # I
}
